$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Shield Parts")
$ws3 = $wb.Worksheets.Item("To Add")

# To Add: E8 text change (new shared string "Test fan with transistor")
$ws3.Range("E8").Value = "Test fan with transistor"

# Shield Parts: B4 text change (new shared string "J7, J15-J26")
$ws2.Range("B4").Value = "J7, J15-J26"
$ws2.Range("F4").Value = 13

# To Add: A4 new text (new shared string "Switch reverse mosfet with schottky")
$ws3.Range("A4").Value = "Switch reverse mosfet with schottky"

$ws2.Range("D4").Select()

# Make "To Add" the active sheet (activeTab=2, tabSelected moves there)
$ws3.Activate()
$ws3.Range("B9").Select()
